$d = $word.ActiveDocument
$d.Content.Find.Execute("61+31=", $true, $true, $false, $false, $false, $true, 1, $false, "78-14=", 2) | Out-Null
$d.Content.Find.Execute("12+67=", $true, $true, $false, $false, $false, $true, 1, $false, "89-33=", 2) | Out-Null
$d.Content.Find.Execute("55+6=", $true, $true, $false, $false, $false, $true, 1, $false, "89-21=", 2) | Out-Null
$d.Content.Find.Execute("44-25=", $true, $true, $false, $false, $false, $true, 1, $false, "73-68=", 2) | Out-Null
$d.Content.Find.Execute("67+32=", $true, $true, $false, $false, $false, $true, 1, $false, "43-21=", 2) | Out-Null
$d.Content.Find.Execute("58-41=", $true, $true, $false, $false, $false, $true, 1, $false, "57+19=", 2) | Out-Null
$d.Content.Find.Execute("17+10=", $true, $true, $false, $false, $false, $true, 1, $false, "76-66=", 2) | Out-Null
$d.Content.Find.Execute("97-74=", $true, $true, $false, $false, $false, $true, 1, $false, "58+30=", 2) | Out-Null
$d.Content.Find.Execute("2+71=", $true, $true, $false, $false, $false, $true, 1, $false, "40+51=", 2) | Out-Null
$d.Content.Find.Execute("88+0=", $true, $true, $false, $false, $false, $true, 1, $false, "16+71=", 2) | Out-Null
$d.Content.Find.Execute("71-47=", $true, $true, $false, $false, $false, $true, 1, $false, "20-9=", 2) | Out-Null
$d.Content.Find.Execute("62-58=", $true, $true, $false, $false, $false, $true, 1, $false, "35+21=", 2) | Out-Null
$d.Content.Find.Execute("55+39=", $true, $true, $false, $false, $false, $true, 1, $false, "12+3=", 2) | Out-Null
$d.Content.Find.Execute("89-40=", $true, $true, $false, $false, $false, $true, 1, $false, "88-8=", 2) | Out-Null
$d.Content.Find.Execute("11-0=", $true, $true, $false, $false, $false, $true, 1, $false, "46+13=", 2) | Out-Null
$d.Content.Find.Execute("45+47=", $true, $true, $false, $false, $false, $true, 1, $false, "75-25=", 2) | Out-Null
$d.Content.Find.Execute("9+8=", $true, $true, $false, $false, $false, $true, 1, $false, "97-59=", 2) | Out-Null
$d.Content.Find.Execute("7+64=", $true, $true, $false, $false, $false, $true, 1, $false, "45+31=", 2) | Out-Null
$d.Content.Find.Execute("5+69=", $true, $true, $false, $false, $false, $true, 1, $false, "13+9=", 2) | Out-Null
$d.Content.Find.Execute("13+48=", $true, $true, $false, $false, $false, $true, 1, $false, "57-19=", 2) | Out-Null
$d.Content.Find.Execute("83-41=", $true, $true, $false, $false, $false, $true, 1, $false, "25+52=", 2) | Out-Null
$d.Content.Find.Execute("73-60=", $true, $true, $false, $false, $false, $true, 1, $false, "37-14=", 2) | Out-Null
$d.Content.Find.Execute("72-12=", $true, $true, $false, $false, $false, $true, 1, $false, "11+60=", 2) | Out-Null
$d.Content.Find.Execute("26+16=", $true, $true, $false, $false, $false, $true, 1, $false, "63-30=", 2) | Out-Null
$d.Content.Find.Execute("4+19=", $true, $true, $false, $false, $false, $true, 1, $false, "59+27=", 2) | Out-Null
$d.Content.Find.Execute("50-8=", $true, $true, $false, $false, $false, $true, 1, $false, "90-87=", 2) | Out-Null
$d.Content.Find.Execute("60+38=", $true, $true, $false, $false, $false, $true, 1, $false, "17+78=", 2) | Out-Null
$d.Content.Find.Execute("97-81=", $true, $true, $false, $false, $false, $true, 1, $false, "64+18=", 2) | Out-Null
$d.Content.Find.Execute("31+2=", $true, $true, $false, $false, $false, $true, 1, $false, "72-18=", 2) | Out-Null
$d.Content.Find.Execute("94-36=", $true, $true, $false, $false, $false, $true, 1, $false, "11+32=", 2) | Out-Null
$d.Content.Find.Execute("98-94=", $true, $true, $false, $false, $false, $true, 1, $false, "12+11=", 2) | Out-Null
$d.Content.Find.Execute("29-17=", $true, $true, $false, $false, $false, $true, 1, $false, "99-75=", 2) | Out-Null
$d.Content.Find.Execute("68-60=", $true, $true, $false, $false, $false, $true, 1, $false, "61-24=", 2) | Out-Null
$d.Content.Find.Execute("52+15=", $true, $true, $false, $false, $false, $true, 1, $false, "70-54=", 2) | Out-Null
$d.Content.Find.Execute("17-11=", $true, $true, $false, $false, $false, $true, 1, $false, "60+5=", 2) | Out-Null
$d.Content.Find.Execute("74-39=", $true, $true, $false, $false, $false, $true, 1, $false, "52+44=", 2) | Out-Null
$d.Content.Find.Execute("64-23=", $true, $true, $false, $false, $false, $true, 1, $false, "38-26=", 2) | Out-Null
$d.Content.Find.Execute("20+25=", $true, $true, $false, $false, $false, $true, 1, $false, "30+53=", 2) | Out-Null
$d.Content.Find.Execute("3+37=", $true, $true, $false, $false, $false, $true, 1, $false, "1+25=", 2) | Out-Null
$d.Content.Find.Execute("4+73=", $true, $true, $false, $false, $false, $true, 1, $false, "8+67=", 2) | Out-Null
$d.Content.Find.Execute("83-38=", $true, $true, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("72-39=", $true, $true, $false, $false, $false, $true, 1, $false, "76-69=", 2) | Out-Null
$d.Content.Find.Execute("84-2=", $true, $true, $false, $false, $false, $true, 1, $false, "75-50=", 2) | Out-Null
$d.Content.Find.Execute("5+10=", $true, $true, $false, $false, $false, $true, 1, $false, "90-61=", 2) | Out-Null
$d.Content.Find.Execute("16+66=", $true, $true, $false, $false, $false, $true, 1, $false, "52-49=", 2) | Out-Null
$d.Content.Find.Execute("39+30=", $true, $true, $false, $false, $false, $true, 1, $false, "65+16=", 2) | Out-Null
$d.Content.Find.Execute("17+8=", $true, $true, $false, $false, $false, $true, 1, $false, "64+26=", 2) | Out-Null
$d.Content.Find.Execute("52+22=", $true, $true, $false, $false, $false, $true, 1, $false, "48+37=", 2) | Out-Null
$d.Content.Find.Execute("94-18=", $true, $true, $false, $false, $false, $true, 1, $false, "29+33=", 2) | Out-Null
$d.Content.Find.Execute("29+27=", $true, $true, $false, $false, $false, $true, 1, $false, "74+14=", 2) | Out-Null
$d.Content.Find.Execute("96-70=", $true, $true, $false, $false, $false, $true, 1, $false, "38+41=", 2) | Out-Null
$d.Content.Find.Execute("1+13=", $true, $true, $false, $false, $false, $true, 1, $false, "58-50=", 2) | Out-Null
$d.Content.Find.Execute("83-65=", $true, $true, $false, $false, $false, $true, 1, $false, "1+56=", 2) | Out-Null
$d.Content.Find.Execute("91-0=", $true, $true, $false, $false, $false, $true, 1, $false, "8+50=", 2) | Out-Null
$d.Content.Find.Execute("26-1=", $true, $true, $false, $false, $false, $true, 1, $false, "65-39=", 2) | Out-Null
$d.Content.Find.Execute("64-58=", $true, $true, $false, $false, $false, $true, 1, $false, "69-4=", 2) | Out-Null
$d.Content.Find.Execute("86-53=", $true, $true, $false, $false, $false, $true, 1, $false, "31+52=", 2) | Out-Null
$d.Content.Find.Execute("40+0=", $true, $true, $false, $false, $false, $true, 1, $false, "69-51=", 2) | Out-Null
$d.Content.Find.Execute("39+43=", $true, $true, $false, $false, $false, $true, 1, $false, "93-44=", 2) | Out-Null
$d.Content.Find.Execute("88-33=", $true, $true, $false, $false, $false, $true, 1, $false, "42-24=", 2) | Out-Null
$d.Content.Find.Execute("30+52=", $true, $true, $false, $false, $false, $true, 1, $false, "37+50=", 2) | Out-Null
$d.Content.Find.Execute("7+33=", $true, $true, $false, $false, $false, $true, 1, $false, "42+44=", 2) | Out-Null
$d.Content.Find.Execute("67-25=", $true, $true, $false, $false, $false, $true, 1, $false, "59+6=", 2) | Out-Null
$d.Content.Find.Execute("2+58=", $true, $true, $false, $false, $false, $true, 1, $false, "51-4=", 2) | Out-Null
$d.Content.Find.Execute("16-13=", $true, $true, $false, $false, $false, $true, 1, $false, "23+64=", 2) | Out-Null
$d.Content.Find.Execute("77+8=", $true, $true, $false, $false, $false, $true, 1, $false, "49+37=", 2) | Out-Null
$d.Content.Find.Execute("32+33=", $true, $true, $false, $false, $false, $true, 1, $false, "54-39=", 2) | Out-Null
$d.Content.Find.Execute("39+9=", $true, $true, $false, $false, $false, $true, 1, $false, "88-19=", 2) | Out-Null
$d.Content.Find.Execute("5+62=", $true, $true, $false, $false, $false, $true, 1, $false, "56-20=", 2) | Out-Null
$d.Content.Find.Execute("68-5=", $true, $true, $false, $false, $false, $true, 1, $false, "20+48=", 2) | Out-Null
$d.Content.Find.Execute("41-33=", $true, $true, $false, $false, $false, $true, 1, $false, "30+47=", 2) | Out-Null
$d.Content.Find.Execute("99-66=", $true, $true, $false, $false, $false, $true, 1, $false, "16+82=", 2) | Out-Null
$d.Content.Find.Execute("79+8=", $true, $true, $false, $false, $false, $true, 1, $false, "26+56=", 2) | Out-Null
$d.Content.Find.Execute("32+60=", $true, $true, $false, $false, $false, $true, 1, $false, "19-3=", 2) | Out-Null
$d.Content.Find.Execute("69+25=", $true, $true, $false, $false, $false, $true, 1, $false, "99-20=", 2) | Out-Null
$d.Content.Find.Execute("50+22=", $true, $true, $false, $false, $false, $true, 1, $false, "74-11=", 2) | Out-Null
$d.Content.Find.Execute("57+6=", $true, $true, $false, $false, $false, $true, 1, $false, "12+25=", 2) | Out-Null
$d.Content.Find.Execute("77-55=", $true, $true, $false, $false, $false, $true, 1, $false, "81-78=", 2) | Out-Null
$d.Content.Find.Execute("48-43=", $true, $true, $false, $false, $false, $true, 1, $false, "39+47=", 2) | Out-Null
$d.Content.Find.Execute("13+70=", $true, $true, $false, $false, $false, $true, 1, $false, "8+0=", 2) | Out-Null
$d.Content.Find.Execute("21-11=", $true, $true, $false, $false, $false, $true, 1, $false, "83-1=", 2) | Out-Null
$d.Content.Find.Execute("94-8=", $true, $true, $false, $false, $false, $true, 1, $false, "23+1=", 2) | Out-Null
$d.Content.Find.Execute("64-0=", $true, $true, $false, $false, $false, $true, 1, $false, "1+60=", 2) | Out-Null
$d.Content.Find.Execute("69+23=", $true, $true, $false, $false, $false, $true, 1, $false, "95-88=", 2) | Out-Null
$d.Content.Find.Execute("67-66=", $true, $true, $false, $false, $false, $true, 1, $false, "19-11=", 2) | Out-Null
$d.Content.Find.Execute("99-55=", $true, $true, $false, $false, $false, $true, 1, $false, "20+47=", 2) | Out-Null
$d.Content.Find.Execute("23-12=", $true, $true, $false, $false, $false, $true, 1, $false, "78+13=", 2) | Out-Null
$d.Content.Find.Execute("1+79=", $true, $true, $false, $false, $false, $true, 1, $false, "27+12=", 2) | Out-Null
$d.Content.Find.Execute("99-91=", $true, $true, $false, $false, $false, $true, 1, $false, "67-22=", 2) | Out-Null
$d.Content.Find.Execute("72-55=", $true, $true, $false, $false, $false, $true, 1, $false, "87-45=", 2) | Out-Null
$d.Content.Find.Execute("53-49=", $true, $true, $false, $false, $false, $true, 1, $false, "11+69=", 2) | Out-Null
$d.Content.Find.Execute("75-44=", $true, $true, $false, $false, $false, $true, 1, $false, "77-32=", 2) | Out-Null
$d.Content.Find.Execute("82+10=", $true, $true, $false, $false, $false, $true, 1, $false, "61+34=", 2) | Out-Null
$d.Content.Find.Execute("9+81=", $true, $true, $false, $false, $false, $true, 1, $false, "17-4=", 2) | Out-Null
$d.Content.Find.Execute("60-48=", $true, $true, $false, $false, $false, $true, 1, $false, "5+52=", 2) | Out-Null
$d.Content.Find.Execute("80-35=", $true, $true, $false, $false, $false, $true, 1, $false, "55-45=", 2) | Out-Null
$d.Content.Find.Execute("64-22=", $true, $true, $false, $false, $false, $true, 1, $false, "70-56=", 2) | Out-Null
$d.Content.Find.Execute("64-18=", $true, $true, $false, $false, $false, $true, 1, $false, "51+45=", 2) | Out-Null
$d.Content.Find.Execute("7+54=", $true, $true, $false, $false, $false, $true, 1, $false, "8+89=", 2) | Out-Null
$d.Content.Find.Execute("76-1=", $true, $true, $false, $false, $false, $true, 1, $false, "34+16=", 2) | Out-Null
